$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("T3").Value = "TrueString"
$ws.Range("U3").Value = 'How the Boolean value True is to be represented in the file. Optional, defaulting to "True".'
$ws.Range("V3").Value = "FalseString"
$ws.Range("W3").Value = 'How the Boolean value False is to be represented in the file. Optional, defaulting to "False".'
$ws.Columns.Item(21).AutoFit()
$ws.Columns.Item(22).AutoFit()
$ws.Columns.Item(23).AutoFit()
